# Apply scheduled market-data refresh to the Leve profit calculations
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row 13
$ws_ALC.Range("H13").Value = 33000
$ws_ALC.Range("J13").Value = 33000
$ws_ALC.Range("L13").Value = 33000
$ws_ALC.Range("N13").Value = -33338

# ALC!row 98
$ws_ALC.Range("H98").Value = 2716.1562
$ws_ALC.Range("I98").Value = 1478.4073
$ws_ALC.Range("J98").Value = 9400
$ws_ALC.Range("K98").Value = 1478.4073
$ws_ALC.Range("L98").Value = 9400
$ws_ALC.Range("M98").Value = 19.59269999999992
$ws_ALC.Range("N98").Value = -12396

# ALC!row 113
$ws_ALC.Range("H113").Value = 13858
$ws_ALC.Range("J113").Value = 13858
$ws_ALC.Range("L113").Value = 13858
$ws_ALC.Range("N113").Value = -20366

# ALC!row 116
$ws_ALC.Range("H116").Value = 1122945.5
$ws_ALC.Range("I116").Value = 3334829.2
$ws_ALC.Range("K116").Value = 3334829.2
$ws_ALC.Range("M116").Value = -3331387.2

# ALC!row 122
$ws_ALC.Range("H122").Value = 2716.1562
$ws_ALC.Range("I122").Value = 1478.4073
$ws_ALC.Range("J122").Value = 9400
$ws_ALC.Range("K122").Value = 4435.2219
$ws_ALC.Range("L122").Value = 28200
$ws_ALC.Range("M122").Value = -1985.2219
$ws_ALC.Range("N122").Value = -33100

# ALC!row 126
$ws_ALC.Range("H126").Value = 42000
$ws_ALC.Range("J126").Value = 42000
$ws_ALC.Range("L126").Value = 42000
$ws_ALC.Range("N126").Value = -51880

# ALC!row 138
$ws_ALC.Range("H138").Value = 2458.04
$ws_ALC.Range("I138").Value = 1167.4783
$ws_ALC.Range("J138").Value = 2843.5325
$ws_ALC.Range("K138").Value = 3502.4349
$ws_ALC.Range("L138").Value = 8530.5975
$ws_ALC.Range("M138").Value = 1637.5651
$ws_ALC.Range("N138").Value = -18810.5975

# ARM!row 45
$ws_ARM.Range("H45").Value = 1623.375
$ws_ARM.Range("I45").Value = 1156.2858
$ws_ARM.Range("K45").Value = 1156.2858
$ws_ARM.Range("M45").Value = -779.2858000000001

# ARM!row 74
$ws_ARM.Range("H74").Value = 2736
$ws_ARM.Range("I74").Value = 2284.8462
$ws_ARM.Range("J74").Value = 3387.6667
$ws_ARM.Range("K74").Value = 2284.8462
$ws_ARM.Range("L74").Value = 3387.6667
$ws_ARM.Range("M74").Value = -1410.8462
$ws_ARM.Range("N74").Value = -5135.6667

# ARM!row 77
$ws_ARM.Range("H77").Value = 2736
$ws_ARM.Range("I77").Value = 2284.8462
$ws_ARM.Range("J77").Value = 3387.6667
$ws_ARM.Range("K77").Value = 11424.231
$ws_ARM.Range("L77").Value = 16938.3335
$ws_ARM.Range("M77").Value = -7056.231
$ws_ARM.Range("N77").Value = -25674.3335

# ARM!row 118
$ws_ARM.Range("H118").Value = 28490
$ws_ARM.Range("J118").Value = 28490
$ws_ARM.Range("L118").Value = 28490
$ws_ARM.Range("N118").Value = -31804

# BSM!row 99
$ws_BSM.Range("H99").Value = 3089.12
$ws_BSM.Range("I99").Value = 1459.2858
$ws_BSM.Range("J99").Value = 5163.4546
$ws_BSM.Range("K99").Value = 1459.2858
$ws_BSM.Range("L99").Value = 5163.4546
$ws_BSM.Range("M99").Value = 38.71419999999989
$ws_BSM.Range("N99").Value = -8159.4546

# BSM!row 125
$ws_BSM.Range("H125").Value = 41776.25
$ws_BSM.Range("J125").Value = 41776.25
$ws_BSM.Range("L125").Value = 41776.25
$ws_BSM.Range("N125").Value = -51616.25

# CRP!row 31
$ws_CRP.Range("H31").Value = 12197972
$ws_CRP.Range("I31").Value = 1141.96
$ws_CRP.Range("K31").Value = 1141.96
$ws_CRP.Range("M31").Value = -846.96

# CRP!row 34
$ws_CRP.Range("H34").Value = 12197972
$ws_CRP.Range("I34").Value = 1141.96
$ws_CRP.Range("K34").Value = 1141.96
$ws_CRP.Range("M34").Value = -939.96

# CRP!row 58
$ws_CRP.Range("H58").Value = 1810.829
$ws_CRP.Range("I58").Value = 1626
$ws_CRP.Range("J58").Value = 2706.5386
$ws_CRP.Range("K58").Value = 1626
$ws_CRP.Range("L58").Value = 2706.5386
$ws_CRP.Range("M58").Value = -1423
$ws_CRP.Range("N58").Value = -3112.5386

# CRP!row 124
$ws_CRP.Range("H124").Value = 30000
$ws_CRP.Range("J124").Value = 30000
$ws_CRP.Range("L124").Value = 30000
$ws_CRP.Range("N124").Value = -34910

# CRP!row 132
$ws_CRP.Range("H132").Value = 3626.262
$ws_CRP.Range("I132").Value = 2887.3333
$ws_CRP.Range("J132").Value = 4365.1904
$ws_CRP.Range("K132").Value = 8661.999899999999
$ws_CRP.Range("L132").Value = 13095.5712
$ws_CRP.Range("M132").Value = -6131.999899999999
$ws_CRP.Range("N132").Value = -18155.5712

# CRP!row 134
$ws_CRP.Range("H134").Value = 6013.52
$ws_CRP.Range("I134").Value = 6921.8823
$ws_CRP.Range("J134").Value = 4083.25
$ws_CRP.Range("K134").Value = 20765.6469
$ws_CRP.Range("L134").Value = 12249.75
$ws_CRP.Range("M134").Value = -18230.6469
$ws_CRP.Range("N134").Value = -17319.75

# CRP!row 136
$ws_CRP.Range("H136").Value = 1810.829
$ws_CRP.Range("I136").Value = 1626
$ws_CRP.Range("J136").Value = 2706.5386
$ws_CRP.Range("K136").Value = 4878
$ws_CRP.Range("L136").Value = 8119.6158
$ws_CRP.Range("M136").Value = -2328
$ws_CRP.Range("N136").Value = -13219.6158

# CRP!row 141
$ws_CRP.Range("H141").Value = 15500
$ws_CRP.Range("J141").Value = 15500
$ws_CRP.Range("L141").Value = 15500
$ws_CRP.Range("N141").Value = -25860

# CUL!row 113
$ws_CUL.Range("H113").Value = 528.8889
$ws_CUL.Range("I113").Value = 517.1579
$ws_CUL.Range("K113").Value = 1551.4737
$ws_CUL.Range("M113").Value = 618.5263

# GSM!row 52
$ws_GSM.Range("H52").Value = 40000
$ws_GSM.Range("J52").Value = 40000
$ws_GSM.Range("L52").Value = 40000
$ws_GSM.Range("N52").Value = -40518

# GSM!row 117
$ws_GSM.Range("H117").Value = 26982
$ws_GSM.Range("J117").Value = 26982
$ws_GSM.Range("L117").Value = 26982
$ws_GSM.Range("N117").Value = -33866

# GSM!row 127
$ws_GSM.Range("H127").Value = 0
$ws_GSM.Range("J127").Value = 0
$ws_GSM.Range("L127").Value = 0
$ws_GSM.Range("N127").ClearContents()

# LTW!row 59
$ws_LTW.Range("H59").Value = 14899
$ws_LTW.Range("J59").Value = 14899
$ws_LTW.Range("L59").Value = 14899
$ws_LTW.Range("N59").Value = -16207

# LTW!row 127
$ws_LTW.Range("H127").Value = 27031.766
$ws_LTW.Range("J127").Value = 27031.766
$ws_LTW.Range("L127").Value = 27031.766
$ws_LTW.Range("N127").Value = -36951.766

# LTW!row 132
$ws_LTW.Range("H132").Value = 3712.7302
$ws_LTW.Range("I132").Value = 2572.2092
$ws_LTW.Range("K132").Value = 7716.6276
$ws_LTW.Range("M132").Value = -5186.6276

# WVR!row 42
$ws_WVR.Range("H42").Value = 50000
$ws_WVR.Range("I42").Value = 0
$ws_WVR.Range("K42").Value = 0
$ws_WVR.Range("M42").ClearContents()

# WVR!row 114
$ws_WVR.Range("H114").Value = 39800
$ws_WVR.Range("J114").Value = 39800
$ws_WVR.Range("L114").Value = 39800
$ws_WVR.Range("N114").Value = -48478

# WVR!row 132
$ws_WVR.Range("H132").Value = 11113187
$ws_WVR.Range("I132").Value = 704.5333000000001
$ws_WVR.Range("J132").Value = 22225670
$ws_WVR.Range("K132").Value = 2113.5999
$ws_WVR.Range("L132").Value = 66677010
$ws_WVR.Range("M132").Value = 416.4000999999998
$ws_WVR.Range("N132").Value = -66682070

# WVR!row 136
$ws_WVR.Range("H136").Value = 1557.9667
$ws_WVR.Range("I136").Value = 920.5789
$ws_WVR.Range("J136").Value = 2658.9092
$ws_WVR.Range("K136").Value = 2761.7367
$ws_WVR.Range("L136").Value = 7976.7276
$ws_WVR.Range("M136").Value = -211.7366999999999
$ws_WVR.Range("N136").Value = -13076.7276
